$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add poll source URL in C1 ---
$ws.Range("C1").Value = "https://leger360.com/wp-content/uploads/2021/01/Legers-North-American-Tracker-January-4th-2021-min.pdf?x16723"

# --- Row 3 header: new column H = "CAN" ---
$ws.Range("H3").Value = "CAN"

# --- New column H values (national numbers) for rows 4-9 ---
$ws.Range("H4").Value = 33
$ws.Range("H5").Value = 33
$ws.Range("H6").Value = 18
$ws.Range("H7").Value = 8
$ws.Range("H8").Value = 6
$ws.Range("H9").Value = 2

# --- New column H values for sample-size rows 10-11 ---
$ws.Range("H10").Value = 1506
$ws.Range("H11").Value = 1506

# --- New row 13: votesum totals with SUM formulas ---
$ws.Range("A13").Value = "votesum"
$ws.Range("B13").Formula = "=SUM(B4:B9)"
$ws.Range("C13").Formula = "=SUM(C4:C9)"
$ws.Range("D13").Formula = "=SUM(D4:D9)"
$ws.Range("E13").Formula = "=SUM(E4:E9)"
$ws.Range("F13").Formula = "=SUM(F4:F9)"
$ws.Range("G13").Formula = "=SUM(G4:G9)"
$ws.Range("H13").Formula = "=SUM(H4:H9)"

# --- Match the surrounding style (style index "1") for cells that changed
#     from the default style to the shared body style in the source edit ---
$bodyStyleSource = $ws.Range("B4")
$ws.Range("A1").Style = $bodyStyleSource.Style
$ws.Range("A2").Style = $bodyStyleSource.Style
$ws.Range("H3:H9").Style = $bodyStyleSource.Style
$ws.Range("A10:H11").Style = $bodyStyleSource.Style
$ws.Range("A13:H13").Style = $bodyStyleSource.Style

# --- Update selection to reflect where the cursor ended up after editing ---
$ws.Range("C14").Select() | Out-Null
